# Applies the "Updated cryptos list" data refresh described by the diff.
# Price (col D) values are textual (dotted thousands separators, etc.), so the
# cell is forced to Text format before assignment to stop Excel from silently
# re-interpreting strings such as "0.999" or "1.00" as numbers; the format is
# then restored to General (matching the original workbook) once the text value
# is locked in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.169.06"
$ws.Range("E2").Value = "  -3.73%  "
$ws.Range("D3").Value = "3.519.30"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "578.16"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").Value = "171.22"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").Value = "3.506.82"
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "0.190"
$ws.Range("E10").Value = "  -6.50%  "
$ws.Range("D11").Value = "6.71"
$ws.Range("E11").Value = "  +13.50%  "
$ws.Range("D12").Value = "0.599"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "47.23"
$ws.Range("E13").Value = "  -4.96%  "
$ws.Range("D14").Value = "0.0000276"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").Value = "685.71"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "4.075.48"
$ws.Range("E16").Value = "  -3.43%  "
$ws.Range("D17").Value = "8.78"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").Value = "69.061.92"
$ws.Range("E18").Value = "  -3.86%  "
$ws.Range("D19").Value = "3.514.57"
$ws.Range("E19").Value = "  -3.28%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "17.41"
$ws.Range("E21").Value = "  -4.93%  "
$ws.Range("D22").Value = "11.13"
$ws.Range("E22").Value = "  -4.14%  "
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("D24").Value = "16.57"
$ws.Range("E24").Value = "  -6.90%  "
$ws.Range("D25").Value = "97.62"
$ws.Range("E25").Value = "  -5.55%  "
$ws.Range("D26").Value = "3.83"
$ws.Range("E26").Value = "  -4.49%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "2.66"
$ws.Range("E28").Value = "  -6.53%  "
$ws.Range("D29").Value = "9.42"
$ws.Range("E29").Value = "  -5.74%  "
$ws.Range("D30").Value = "33.24"
$ws.Range("E30").Value = "  -5.62%  "
$ws.Range("D31").Value = "8.84"
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("D32").Value = "3.18"
$ws.Range("E32").Value = "  -6.96%  "
$ws.Range("E33").Value = "  -5.83%  "
$ws.Range("D34").Value = "7.24"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").Value = "569.75"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("D36").Value = "3.66"
$ws.Range("E36").Value = "  -12.82%  "
$ws.Range("D37").Value = "10.84"
$ws.Range("E37").Value = "  -4.13%  "
$ws.Range("E38").Value = "  -3.35%  "
$ws.Range("D39").Value = "57.17"
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "0.138"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0440"
$ws.Range("E42").Value = "  -6.03%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "0.337"
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.440.91"
$ws.Range("E44").Value = "  -6.40%  "
$ws.Range("D45").Value = "33.26"
$ws.Range("E45").Value = "  -6.74%  "
$ws.Range("D46").Value = "0.0₃0704"
$ws.Range("E46").Value = "  -7.88%  "
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("D48").Value = "2.59"
$ws.Range("E48").Value = "  -7.02%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "134.47"
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("D51").Value = "0.149"
$ws.Range("E51").Value = "  -0.71%  "

$ws.Range("D2:D51").Style = "Normal"

